# Updated cryptos list on Tue Jun  4 08:59:41 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.909.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.762.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "627.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.758.76"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("E10").Value = "  -2.62%  "
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.87%  "
$ws.Range("E13").Value = "  -5.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.406.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.798.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.949.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("E18").Value = "  -2.76%  "
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "467.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.701"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("E25").Value = "  -7.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.912.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.27%  "
$ws.Range("E34").Value = "  +18.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "28.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.97%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.718.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.41%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.88%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.960"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "156.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.76%  "
